$d = $word.ActiveDocument

# The first paragraph carries the "Title" style and is currently empty.
# We need to add two runs to it:
#   1. "<blockTable>" in Courier font
#   2. " Tag Demo" in the default (no explicit) formatting
$titlePara = $d.Paragraphs(1)

$tagText = "<blockTable>"
$restText = " Tag Demo"

# Insert all the new text at the very start of the title paragraph
# (i.e. before its paragraph mark).
$titlePara.Range.InsertBefore($tagText + $restText)

# Re-fetch positions: the paragraph now starts at its original start
# position and the inserted text occupies the first
# Len(tagText + restText) characters of it.
$paraStart = $titlePara.Range.Start

# Format just the "<blockTable>" portion with the Courier font.
$tagRange = $d.Range($paraStart, $paraStart + $tagText.Length)
$tagRange.Font.Name = "Courier"
